# Generate Report for Handback
# The 694eb7ed-... file failed to hand back: the generated handback file
# name did not match the expected handoff file name. Update the Status
# cells for that row on every sheet and add an Error Detail message on the
# per-language sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhErrorDetail = "Handback file name: u0yyxug1.3p0 is different with handoff file name: 694eb7ed-d2b8-41f4-a4bd-d04ef96d3ec6.09a500c58d892316bf5d36ee3faf2320bbbaec1b.zh-cn."
$deErrorDetail = "Handback file name: u0yyxug1.3p0 is different with handoff file name: 694eb7ed-d2b8-41f4-a4bd-d04ef96d3ec6.09a500c58d892316bf5d36ee3faf2320bbbaec1b.de-de."

# Overview sheet: row 3 is the 694eb7ed-...md file. Column B = zh-cn status,
# column C = de-de status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# zh-cn sheet: row 3 is the 694eb7ed-...md file. Column C = Status,
# column K = Error Detail.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("K3").Value = $zhErrorDetail

# de-de sheet: row 3 is the 694eb7ed-...md file. Column C = Status,
# column K = Error Detail.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("K3").Value = $deErrorDetail
